{"js": "// The template's item table uses bracketed placeholders like [xcode],\n// [barCode], [unitName], etc. The \"quantity\" column placeholder was\n// renamed from [count] to [qty] (table column rename erp -> srm).\n//\n// Replace the \"count\" placeholder token with \"qty\", leaving the\n// surrounding square brackets untouched.\nconst body = context.document.body;\n\nconst results = body.search(\"count\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"qty\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The template's item table uses bracketed placeholders like [xcode],\n# [barCode], [unitName], etc. The \"quantity\" column placeholder was\n# renamed from [count] to [qty] (table column rename erp -> srm).\n#\n# Replace the \"count\" placeholder token with \"qty\", leaving the\n# surrounding square brackets untouched.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"count\"\n$find.Replacement.Text = \"qty\"\n$find.Forward = $true\n$find.Wrap = 1            # wdFindContinue\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n# wdReplaceAll = 2\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n"}
